# Apply the quarterly database update + column shift described by the commit
# 'update database and change read_price algorithm'.
#
# The sheet holds 10 reporting-period columns (D:M). Each quarterly refresh
# drops the oldest period (old column D), shifts the remaining 9 periods one
# column to the left, and appends the newest period's figures into column M.
# The period/date header labels in row 8 and row 9 follow the same left-shift,
# with a brand new period label/date appended at M8/M9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: 'دوره مالی' (financial period) header labels, shifted left one column ---
$row8 = @{ "D" = "12 ماهه منتهی به 1399/08"; "E" = "3 ماهه منتهی به 1399/11"; "F" = "6 ماهه منتهی به 1400/02"; "G" = "9 ماهه منتهی به 1400/05"; "H" = "12 ماهه منتهی به 1400/08"; "I" = "3 ماهه منتهی به 1400/11"; "J" = "6 ماهه منتهی به 1401/02"; "K" = "9 ماهه منتهی به 1401/05"; "L" = "12 ماهه منتهی به 1401/08"; "M" = "3 ماهه منتهی به 1401/11" }
foreach ($col in $row8.Keys) {
    $ws.Range("$col" + "8").Value = $row8[$col]
}

# --- Row 9: 'تاریخ انتشار' (publish date) header labels, shifted left one column ---
$row9 = @{ "D" = "1400-12-09 (9)"; "E" = "1400-12-28 (2)"; "F" = "1401-05-12 (4)"; "G" = "1401-06-30 (2)"; "H" = "1401-12-13 (9)"; "I" = "1401-12-28 (2)"; "J" = "1401-05-12 (2)"; "K" = "1401-06-30"; "L" = "1401-12-28 (3)"; "M" = "1401-12-28" }
foreach ($col in $row9.Keys) {
    $ws.Range("$col" + "9").Value = $row9[$col]
}

# --- Data rows 11-27: financial figures, shifted left one column, new quarter in M ---
$dataRows = @{
    11 = @{ "D" = 5209685; "E" = 2291174; "F" = 5190258; "G" = 7411215; "H" = 10408824; "I" = 2955544; "J" = 6078079; "K" = 9901046; "L" = 14266243; "M" = 3612004 }
    12 = @{ "D" = -3350233; "E" = -1402084; "F" = -3161767; "G" = -4639229; "H" = -7524301; "I" = -1977948; "J" = -4094200; "K" = -6940360; "L" = -10590781; "M" = -2599981 }
    13 = @{ "D" = 1859452; "E" = 889090; "F" = 2028491; "G" = 2771986; "H" = 2884523; "I" = 977596; "J" = 1983879; "K" = 2960686; "L" = 3675462; "M" = 1012023 }
    14 = @{ "D" = -120107; "E" = -28555; "F" = -161847; "G" = -113026; "H" = -276573; "I" = -42272; "J" = -244780; "K" = -260733; "L" = -322707; "M" = -96214 }
    15 = @{ "D" = 0; "E" = 0; "F" = 0; "G" = 0; "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0; "M" = 0 }
    16 = @{ "D" = 26832; "E" = 2204; "F" = -10132; "G" = 111; "H" = -30045; "I" = 0; "J" = 75497; "K" = 95029; "L" = 39809; "M" = 146000 }
    17 = @{ "D" = 1766177; "E" = 862739; "F" = 1856512; "G" = 2659071; "H" = 2577905; "I" = 935324; "J" = 1814596; "K" = 2794982; "L" = 3392564; "M" = 1061809 }
    18 = @{ "D" = -332949; "E" = -145984; "F" = -258811; "G" = -304264; "H" = -515122; "I" = -144270; "J" = -280649; "K" = -478062; "L" = -654109; "M" = -194505 }
    19 = @{ "D" = 47884; "E" = 240; "F" = 36694; "G" = 15969; "H" = 454609; "I" = 0; "J" = 65439; "K" = 91862; "L" = 42988; "M" = 0 }
    20 = @{ "D" = 1481112; "E" = 716995; "F" = 1634395; "G" = 2370776; "H" = 2517392; "I" = 791054; "J" = 1599386; "K" = 2408782; "L" = 2781443; "M" = 867304 }
    21 = @{ "D" = -282112; "E" = -139173; "F" = -374824; "G" = -462798; "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0; "M" = 0 }
    22 = @{ "D" = 1199000; "E" = 577822; "F" = 1259571; "G" = 1907978; "H" = 2517392; "I" = 791054; "J" = 1599386; "K" = 2408782; "L" = 2781443; "M" = 867304 }
    23 = @{ "D" = 0; "E" = 0; "F" = 0; "G" = 0; "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0; "M" = 0 }
    24 = @{ "D" = 1199000; "E" = 577822; "F" = 1259571; "G" = 1907978; "H" = 2517392; "I" = 791054; "J" = 1599386; "K" = 2408782; "L" = 2781443; "M" = 867304 }
    25 = @{ "D" = 128; "E" = 62; "F" = 135; "G" = 204; "H" = 269; "I" = 85; "J" = 171; "K" = 257; "L" = 297; "M" = 93 }
    26 = @{ "D" = 9358065; "E" = 9358065; "F" = 9358065; "G" = 9358065; "H" = 9358065; "I" = 9358065; "J" = 9358065; "K" = 9358065; "L" = 9358065; "M" = 9358065 }
    27 = @{ "D" = 128; "E" = 62; "F" = 135; "G" = 204; "H" = 269; "I" = 85; "J" = 171; "K" = 257; "L" = 297; "M" = 93 }
}
foreach ($r in $dataRows.Keys) {
    $rowData = $dataRows[$r]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col" + "$r").Value = $rowData[$col]
    }
}

